$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36 to make room for the Sugarcane entry.
$ws.Rows(36).Insert()

# Copy the formatting (border/font/alignment) from the row that used to be
# row 36 (now shifted down to row 37) so the new row matches the sheet's
# existing look.
$ws.Range("A37:E37").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match the row height shown for the new Sugarcane row.
$ws.Rows(36).RowHeight = 34.5

# Fill in the new Sugarcane resource entry.
$ws.Range("A36").Value = "SUGARCANE"
$ws.Range("B36").Value = "Known as ""reeds that produce honey without bees. Sugarcane is a very tall plant that can grow up to 20 feet tall. The tall stalks are rich in sugar and can be processed in various ways or eaten raw. "
$ws.Range("C36").Value = "3 coins per stalk"
$ws.Range("D36").Value = "Food, Crop, Plant, Ingredient"
$ws.Range("E36").Value = "A long reed filled with raw sugar."

# Update the saved view state to match what was recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("E36").Select()
